# credData.xlsx - "dynamic object creation (passing dynammic xpath) for
# rsOrderHistoryConfirm" - apply the recorded edits via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Shared string fix-up: the RSURL test-data cell (M2) pointed at
#    "https://mirandakate.cabitest3.com/" - drop the trailing slash.
# ------------------------------------------------------------------
$ws.Range("M2").Value = "https://mirandakate.cabitest3.com"

# ------------------------------------------------------------------
# 2) Column G (the CCURL column) needs to be widened so the longer
#    value shows in full - go from ~31.7 chars to ~48.9 chars.
#    (ColumnWidth is expressed in characters; Excel snaps this to a
#    pixel grid internally, 48 is the closest input that lands on the
#    intended width.)
# ------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 48

# ------------------------------------------------------------------
# 3) Scroll the sheet view so column G is the left-most visible
#    column instead of column D (topLeftCell D1 -> G1).
# ------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 7
$win.ScrollRow = 1

# ------------------------------------------------------------------
# 4) The workbook was re-saved from a different folder on disk, so
#    Excel updated the recorded absolute path of the workbook:
#      C:\Users\user\Documents\testautomation_new\NewArrivals\
#      -> C:\Users\user\Documents\testautomation\NewArrivals\
#    This value (x15ac:absPath) is stamped by Excel itself from the
#    real save location and is not exposed as a settable property on
#    the Workbook/Application object model (there is no FullName /
#    Path setter), so it cannot be changed through COM automation -
#    it is left as-is here.
# ------------------------------------------------------------------

$wb.Save()
